$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text representation
# (values such as "1.70" or "74.512.94" would otherwise be re-interpreted
# as numbers and lose their exact formatting), so force Text format first.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "73.982.99"
$ws.Range("E2").Value = "  +8.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.573.85"
$ws.Range("E3").Value = "  +6.10%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "184.51"
$ws.Range("E5").Value = "  +15.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "576.65"
$ws.Range("E6").Value = "  +3.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("E8").Value = "  +4.75%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.201"
$ws.Range("E9").Value = "  +23.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.571.34"
$ws.Range("E10").Value = "  +6.05%  "

$ws.Range("E11").Value = "  -0.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.358"
$ws.Range("E12").Value = "  +8.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.75"
$ws.Range("E13").Value = "  +3.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000189"
$ws.Range("E14").Value = "  +9.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "74.087.51"
$ws.Range("E15").Value = "  +8.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.046.86"
$ws.Range("E16").Value = "  +6.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.04"
$ws.Range("E17").Value = "  +13.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.577.69"
$ws.Range("E18").Value = "  +6.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.52"
$ws.Range("E19").Value = "  +23.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.64"
$ws.Range("E20").Value = "  +11.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.97"
$ws.Range("E21").Value = "  +12.17%  "

$ws.Range("E22").Value = "  +20.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.04"
$ws.Range("E23").Value = "  +5.73%  "

$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.31"
$ws.Range("E25").Value = "  +4.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.13"
$ws.Range("E26").Value = "  +12.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.17"
$ws.Range("E27").Value = "  +12.01%  "

$ws.Range("E28").Value = "  +6.33%  "

$ws.Range("E29").Value = "  -0.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0932"
$ws.Range("E30").Value = "  +14.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.89"
$ws.Range("E31").Value = "  +10.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "497.67"
$ws.Range("E32").Value = "  +17.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.33"
$ws.Range("E33").Value = "  +16.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.70"
$ws.Range("E34").Value = "  +5.78%  "

$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.118"
$ws.Range("E36").Value = "  +12.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.89"
$ws.Range("E37").Value = "  -0.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.10"
$ws.Range("E38").Value = "  +6.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.35"
$ws.Range("E39").Value = "  +1.63%  "

$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.88"
$ws.Range("E41").Value = "  +13.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.66"
$ws.Range("E42").Value = "  +12.19%  "

$ws.Range("E43").Value = "  +7.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.41"
$ws.Range("E44").Value = "  +19.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "38.91"
$ws.Range("E45").Value = "  +4.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.15"
$ws.Range("E46").Value = "  +7.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.45"
$ws.Range("E47").Value = "  +12.52%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0818"
$ws.Range("E48").Value = "  +14.71%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.59"
$ws.Range("E49").Value = "  +7.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.517"
$ws.Range("E50").Value = "  +7.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0965"
$ws.Range("E51").Value = "  +5.64%  "
